# Add a new review record in row 16 of Sheet1, mirroring the layout of the
# existing rows (e.g. row 15): appid/keyword in A:B, two hyperlinked emails
# in C:D, a timestamp in E, a review comment in F and a "recovery" flag in G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Start from row 15's formatting (number formats / styles / borders) so the
#    new row inherits the same cell styles (s="1" for A, s="2" for C/D, etc.)
#    instead of picking up the workbook's default style.
$ws.Range("A15:G15").Copy()
$ws.Range("A16:G16").PasteSpecial(-4122)   # xlPasteFormats

# 2) Fill in the new row's values.
$ws.Range("A16").Value = "com.hamxa.shaynachim"
$ws.Range("B16").Value = "bitcoin"
$ws.Range("C16").Value = "rontiddler560@gmail.com"
$ws.Range("D16").Value = "halachme@gmail.com"
$ws.Range("E16").Value = "27/5/2019 15:59"
$ws.Range("F16").Value = "money money money.. for beginners only"
$ws.Range("G16").Value = "no"

# 3) Turn the two e-mail addresses into mailto: hyperlinks, same as the other
#    rows in the sheet.
$ws.Hyperlinks.Add($ws.Range("C16"), "mailto:rontiddler560@gmail.com", [Type]::Missing, [Type]::Missing, "rontiddler560@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D16"), "mailto:halachme@gmail.com", [Type]::Missing, [Type]::Missing, "halachme@gmail.com")

# 4) Adding the hyperlinks can reset the C16/D16 cell style to the built-in
#    "Hyperlink" style, so re-apply row 15's formatting on top of the values.
$ws.Range("A15:G15").Copy()
$ws.Range("A16:G16").PasteSpecial(-4122)   # xlPasteFormats

# 5) The new row is a touch shorter than the surrounding rows (12.8pt instead
#    of 13.8pt).
$ws.Rows.Item(16).RowHeight = 12.8

# 6) Reflect the sheet's new selection (the user had just finished entering
#    the new record in C16:D16).
$ws.Range("C16:D16").Select()
